# Update the "Assigné à" column (D) of the Kanban table on sheet "Tableau KANBAN"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tableau KANBAN")

$ws.Range("D6").Value  = "Teddy"
$ws.Range("D7").Value  = " Teddy "
$ws.Range("D8").Value  = " Teddy     "
$ws.Range("D9").Value  = " Ny       "
$ws.Range("D10").Value = " Ny     "
$ws.Range("D11").Value = " Teddy&Ny    "
$ws.Range("D12").Value = "  Teddy&Ny        "

# Update the active cell selection to D12
$ws.Activate()
$ws.Range("D12").Select()
